$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 0.01500074871952251
$ws.Cells.Item(2, 4).Value = 0.1075639941003912
$ws.Cells.Item(2, 5).Value = 0.13266765208677
$ws.Cells.Item(2, 6).Value = 2.729329427611674
$ws.Cells.Item(2, 7).Value = 2.200389605886699
$ws.Cells.Item(2, 8).Value = 1.609776939536573
$ws.Cells.Item(2, 9).Value = 1.648694371904597
$ws.Cells.Item(2, 10).Value = 0.2129241647032103

$ws.Cells.Item(3, 3).Value = 0.01495475488447084
$ws.Cells.Item(3, 4).Value = 0.1077514045506696
$ws.Cells.Item(3, 5).Value = 0.1316056239947621
$ws.Cells.Item(3, 6).Value = 2.614574432376259
$ws.Cells.Item(3, 7).Value = 2.070184954696032
$ws.Cells.Item(3, 8).Value = 1.556616617358628
$ws.Cells.Item(3, 9).Value = 1.565173675784195
$ws.Cells.Item(3, 10).Value = 0.2067062043403638

$ws.Cells.Item(4, 3).Value = 0.01493494521318794
$ws.Cells.Item(4, 4).Value = 0.107875866184276
$ws.Cells.Item(4, 5).Value = 0.1310269187202096
$ws.Cells.Item(4, 6).Value = 2.546056338326565
$ws.Cells.Item(4, 7).Value = 1.991836982531964
$ws.Cells.Item(4, 8).Value = 1.525089552693203
$ws.Cells.Item(4, 9).Value = 1.515082758864821
$ws.Cells.Item(4, 10).Value = 0.2030360327594138

$ws.Cells.Item(5, 3).Value = 0.01492897907347412
$ws.Cells.Item(5, 4).Value = 0.1079289492167297
$ws.Cells.Item(5, 5).Value = 0.13080945999139
$ws.Cells.Item(5, 6).Value = 2.518616852301335
$ws.Cells.Item(5, 7).Value = 1.960304707828357
$ws.Cells.Item(5, 8).Value = 1.512519096264384
$ws.Cells.Item(5, 9).Value = 1.494965473062436
$ws.Cells.Item(5, 10).Value = 0.2015771646160616

$ws.Cells.Item(6, 3).Value = 0.01492811523809934
$ws.Cells.Item(6, 4).Value = 0.10793790646543
$ws.Cells.Item(6, 5).Value = 0.1307744581071724
$ws.Cells.Item(6, 6).Value = 2.514089493075005
$ws.Cells.Item(6, 7).Value = 1.955092470030536
$ws.Cells.Item(6, 8).Value = 1.510448431505552
$ws.Cells.Item(6, 9).Value = 1.4916427149453
$ws.Cells.Item(6, 10).Value = 0.2013371297368494

$ws.Cells.Item(7, 3).Value = 0.01493485623934276
$ws.Cells.Item(7, 4).Value = 0.107876572506173
$ws.Cells.Item(7, 5).Value = 0.1310239117327185
$ws.Cells.Item(7, 6).Value = 2.545684335439148
$ws.Cells.Item(7, 7).Value = 1.991410136056288
$ws.Cells.Item(7, 8).Value = 1.524918905083638
$ws.Cells.Item(7, 9).Value = 1.514810260670487
$ws.Cells.Item(7, 10).Value = 0.2030162096223904

$ws.Cells.Item(8, 3).Value = 0.01498313195366308
$ws.Cells.Item(8, 4).Value = 0.1076266656038589
$ws.Cells.Item(8, 5).Value = 0.1322861816228844
$ws.Cells.Item(8, 6).Value = 2.689355098671854
$ws.Cells.Item(8, 7).Value = 2.155159211585641
$ws.Cells.Item(8, 8).Value = 1.591214328018481
$ws.Cells.Item(8, 9).Value = 1.619646549083043
$ws.Cells.Item(8, 10).Value = 0.2107493281297366

$ws.Cells.Item(9, 3).Value = 0.0151453429096442
$ws.Cells.Item(9, 4).Value = 0.1072110100357313
$ws.Cells.Item(9, 5).Value = 0.1353479623010578
$ws.Cells.Item(9, 6).Value = 2.986798958259925
$ws.Cells.Item(9, 7).Value = 2.489277187964319
$ws.Cells.Item(9, 8).Value = 1.730198639824835
$ws.Cells.Item(9, 9).Value = 1.834892252438209
$ws.Cells.Item(9, 10).Value = 0.2271041262268199

$ws.Cells.Item(10, 3).Value = 0.0153066505152637
$ws.Cells.Item(10, 4).Value = 0.1069508684809621
$ws.Cells.Item(10, 5).Value = 0.1379613603961403
$ws.Cells.Item(10, 6).Value = 3.21534896242585
$ws.Cells.Item(10, 7).Value = 2.743165916485282
$ws.Cells.Item(10, 8).Value = 1.837998461700465
$ws.Cells.Item(10, 9).Value = 1.999237326983689
$ws.Cells.Item(10, 10).Value = 0.2398729209588879

$ws.Cells.Item(11, 3).Value = 0.01538939953681506
$ws.Cells.Item(11, 4).Value = 0.1068423256005708
$ws.Cells.Item(11, 5).Value = 0.1392307378497932
$ws.Cells.Item(11, 6).Value = 3.32159736931149
$ws.Cells.Item(11, 7).Value = 2.860601604875853
$ws.Cells.Item(11, 8).Value = 1.888322784717218
$ws.Cells.Item(11, 9).Value = 2.075419860908085
$ws.Cells.Item(11, 10).Value = 0.2458514306728006

$ws.Cells.Item(12, 3).Value = 0.01542209969176156
$ws.Cells.Item(12, 4).Value = 0.1068026307307428
$ws.Cells.Item(12, 5).Value = 0.1397231131013577
$ws.Cells.Item(12, 6).Value = 3.362166558345194
$ws.Cells.Item(12, 7).Value = 2.905359000080978
$ws.Cells.Item(12, 8).Value = 1.907567967969442
$ws.Cells.Item(12, 9).Value = 2.10447810318459
$ws.Cells.Item(12, 10).Value = 0.2481402578854244

$ws.Cells.Item(13, 3).Value = 0.01541499613716013
$ws.Cells.Item(13, 4).Value = 0.1068111171222066
$ws.Cells.Item(13, 5).Value = 0.1396165496052433
$ws.Cells.Item(13, 6).Value = 3.353414228960958
$ws.Cells.Item(13, 7).Value = 2.895706797964522
$ws.Cells.Item(13, 8).Value = 1.903414733845125
$ws.Cells.Item(13, 9).Value = 2.098210485753555
$ws.Cells.Item(13, 10).Value = 0.2476462033039724

$ws.Cells.Item(14, 3).Value = 0.01539206234078705
$ws.Cells.Item(14, 4).Value = 0.1068390316654844
$ws.Cells.Item(14, 5).Value = 0.1392710109424584
$ws.Cells.Item(14, 6).Value = 3.324928259606082
$ws.Cells.Item(14, 7).Value = 2.864278023235613
$ws.Cells.Item(14, 8).Value = 1.889902299789014
$ws.Cells.Item(14, 9).Value = 2.07780626958774
$ws.Cells.Item(14, 10).Value = 0.2460392322431062

$ws.Cells.Item(15, 3).Value = 0.01537819301030652
$ws.Cells.Item(15, 4).Value = 0.1068563134678397
$ws.Cells.Item(15, 5).Value = 0.1390608841685363
$ws.Cells.Item(15, 6).Value = 3.307523659011053
$ws.Cells.Item(15, 7).Value = 2.845064614609498
$ws.Cells.Item(15, 8).Value = 1.881650193900043
$ws.Cells.Item(15, 9).Value = 2.065335547081503
$ws.Cells.Item(15, 10).Value = 0.2450581720557352

$ws.Cells.Item(16, 3).Value = 0.01530143269788198
$ws.Cells.Item(16, 4).Value = 0.1069581590765658
$ws.Cells.Item(16, 5).Value = 0.1378800333383303
$ws.Cells.Item(16, 6).Value = 3.208451882713092
$ws.Cells.Item(16, 7).Value = 2.735530945545804
$ws.Cells.Item(16, 8).Value = 1.834735823243136
$ws.Cells.Item(16, 9).Value = 1.994287658834168
$ws.Cells.Item(16, 10).Value = 0.23948566966844

$ws.Cells.Item(17, 3).Value = 0.0152567547858844
$ws.Cells.Item(17, 4).Value = 0.1070231466110485
$ws.Cells.Item(17, 5).Value = 0.1371763219990818
$ws.Cells.Item(17, 6).Value = 3.14826367969107
$ws.Cells.Item(17, 7).Value = 2.668838223621094
$ws.Cells.Item(17, 8).Value = 1.806287179186654
$ws.Cells.Item(17, 9).Value = 1.951069736925746
$ws.Cells.Item(17, 10).Value = 0.2361109650743458

$ws.Cells.Item(18, 3).Value = 0.0152319380363366
$ws.Cells.Item(18, 4).Value = 0.1070614478785616
$ws.Cells.Item(18, 5).Value = 0.1367791389186905
$ws.Cells.Item(18, 6).Value = 3.113859156461672
$ws.Cells.Item(18, 7).Value = 2.630660485871545
$ws.Cells.Item(18, 8).Value = 1.790045062812396
$ws.Cells.Item(18, 9).Value = 1.926345376060937
$ws.Cells.Item(18, 10).Value = 0.2341858897510747

$ws.Cells.Item(19, 3).Value = 0.01522368626370607
$ws.Cells.Item(19, 4).Value = 0.1070745744399382
$ws.Cells.Item(19, 5).Value = 0.1366459569529717
$ws.Cells.Item(19, 6).Value = 3.102246932674007
$ws.Cells.Item(19, 7).Value = 2.617765203410215
$ws.Cells.Item(19, 8).Value = 1.784566399661514
$ws.Cells.Item(19, 9).Value = 1.917996894752775
$ws.Cells.Item(19, 10).Value = 0.2335368205613122

$ws.Cells.Item(20, 3).Value = 0.01526141954071036
$ws.Cells.Item(20, 4).Value = 0.1070161331441852
$ws.Cells.Item(20, 5).Value = 0.1372504487062614
$ws.Cells.Item(20, 6).Value = 3.154648601167054
$ws.Cells.Item(20, 7).Value = 2.675918868217252
$ws.Cells.Item(20, 8).Value = 1.80930305883669
$ws.Cells.Item(20, 9).Value = 1.95565651010952
$ws.Cells.Item(20, 10).Value = 0.2364685524690486

$ws.Cells.Item(21, 3).Value = 0.01539876136780194
$ws.Cells.Item(21, 4).Value = 0.1068307942810289
$ws.Cells.Item(21, 5).Value = 0.1393721858826531
$ws.Cells.Item(21, 6).Value = 3.333286121714707
$ws.Cells.Item(21, 7).Value = 2.873501560967327
$ws.Cells.Item(21, 8).Value = 1.893866085626826
$ws.Cells.Item(21, 9).Value = 2.083793752641611
$ws.Cells.Item(21, 10).Value = 0.2465105592145989

$ws.Cells.Item(22, 3).Value = 0.01549648643440449
$ws.Cells.Item(22, 4).Value = 0.1067178711083758
$ws.Cells.Item(22, 5).Value = 0.1408270527599811
$ws.Cells.Item(22, 6).Value = 3.451993562778256
$ws.Cells.Item(22, 7).Value = 3.0043102689649
$ws.Cells.Item(22, 8).Value = 1.950233095680687
$ws.Cells.Item(22, 9).Value = 2.168762892302453
$ws.Cells.Item(22, 10).Value = 0.2532189225224784

$ws.Cells.Item(23, 3).Value = 0.01544359403042606
$ws.Cells.Item(23, 4).Value = 0.1067773896412065
$ws.Cells.Item(23, 5).Value = 0.1400442865526053
$ws.Cells.Item(23, 6).Value = 3.388455569686528
$ws.Cells.Item(23, 7).Value = 2.934339055916894
$ws.Cells.Item(23, 8).Value = 1.920047087610783
$ws.Cells.Item(23, 9).Value = 2.123299507958961
$ws.Cells.Item(23, 10).Value = 0.2496250895907792

$ws.Cells.Item(24, 3).Value = 0.01525930789880192
$ws.Cells.Item(24, 4).Value = 0.1070193010077216
$ws.Cells.Item(24, 5).Value = 0.1372169130121357
$ws.Cells.Item(24, 6).Value = 3.151761360513774
$ws.Cells.Item(24, 7).Value = 2.672717196094936
$ws.Cells.Item(24, 8).Value = 1.807939226738313
$ws.Cells.Item(24, 9).Value = 1.953582449725758
$ws.Cells.Item(24, 10).Value = 0.2363068402489432

$ws.Cells.Item(25, 3).Value = 0.01509413279875815
$ws.Cells.Item(25, 4).Value = 0.1073155030958837
$ws.Cells.Item(25, 5).Value = 0.1344562352988277
$ws.Cells.Item(25, 6).Value = 2.904605856285286
$ws.Cells.Item(25, 7).Value = 2.397447652994629
$ws.Cells.Item(25, 8).Value = 1.691616583404766
$ws.Cells.Item(25, 9).Value = 1.775596163979088
$ws.Cells.Item(25, 10).Value = 0.2225495127181745
